# Fruta / hortaliza, semanal
# Insert a new daily-price record as the first row of the "Plátano" block
# (row 282), pushing all subsequent rows down by one. This grows the used
# range from A1:T342 to A1:T343.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 282; Excel shifts row 282..342 down to 283..343
# and carries the existing row 282 formatting (incl. the date style on column D)
# into the new, now-empty row 282.
$ws.Rows.Item(282).Insert()

# Populate the new row 282 with the new record's data (same market/product
# template as the surrounding rows: Vega Modelo de Temuco / La Araucanía /
# Fruta / Plátano, priced $/caja 20 kilos, origin Ecuador).
$ws.Cells.Item(282, 1).Value = 10
$ws.Cells.Item(282, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(282, 3).Value = "La Araucanía"
$ws.Cells.Item(282, 4).Value = 44476
$ws.Cells.Item(282, 5).Value = 9
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100108
$ws.Cells.Item(282, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(282, 9).Value = 100108006
$ws.Cells.Item(282, 10).Value = "Plátano"
$ws.Cells.Item(282, 11).Value = "Sin especificar"
$ws.Cells.Item(282, 12).Value = "Pintón"
$ws.Cells.Item(282, 13).Value = 2080
$ws.Cells.Item(282, 14).Value = 21000
$ws.Cells.Item(282, 15).Value = 23000
$ws.Cells.Item(282, 16).Value = 22154
$ws.Cells.Item(282, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(282, 18).Value = "Ecuador"
$ws.Cells.Item(282, 19).Value = 1108
$ws.Cells.Item(282, 20).Value = 20
